$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their original text formatting so that
# numeric-looking strings (e.g. "1.008", "29.300.46") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '29.300.46'
$ws.Range('E2').Value = '  -0.17%  '
# Row 3
$ws.Range('D3').Value = '1.868.50'
$ws.Range('E3').Value = '  -0.56%  '
# Row 4
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.76%  '
# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '242.78'
$ws.Range('E5').Value = '  +0.24%  '
# Row 6
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = '0.6978'
$ws.Range('E6').Value = '  -2.26%  '
# Row 7
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.63%  '
# Row 8
$ws.Range('D8').Value = '0.07761'
$ws.Range('E8').Value = '  -3.94%  '
# Row 9
$ws.Range('D9').Value = '0.3094'
$ws.Range('E9').Value = '  -1.12%  '
# Row 10
$ws.Range('D10').Value = '24.05'
$ws.Range('E10').Value = '  -4.74%  '
# Row 11
$ws.Range('D11').Value = '0.08011'
$ws.Range('E11').Value = '  -4.11%  '
# Row 12
$ws.Range('D12').Value = '1.878.29'
$ws.Range('E12').Value = '  +0.23%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.146'
$ws.Range('E13').Value = '  -1.94%  '
# Row 14
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = '92.72'
$ws.Range('E14').Value = '  +1.35%  '
# Row 15
$ws.Range('D15').Value = '0.6922'
$ws.Range('E15').Value = '  -3.73%  '
# Row 16
$ws.Range('D16').Value = '6.347'
$ws.Range('E16').Value = '  +1.27%  '
# Row 17
$ws.Range('D17').Value = '29.364.48'
$ws.Range('E17').Value = '  +0.06%  '
# Row 18
$ws.Range('D18').Value = '0.000008268'
$ws.Range('E18').Value = '  -1.50%  '
# Row 19
$ws.Range('D19').Value = '249.31'
$ws.Range('E19').Value = '  +3.54%  '
# Row 20
$ws.Range('D20').Value = '2.140.11'
$ws.Range('E20').Value = '  +0.70%  '
# Row 21
$ws.Range('D21').Value = '13.11'
$ws.Range('E21').Value = '  -1.03%  '
# Row 22
$ws.Range('D22').Value = '1.005'
$ws.Range('E22').Value = '  +0.45%  '
# Row 23
$ws.Range('D23').Value = '7.510'
$ws.Range('E23').Value = '  -3.69%  '
# Row 24
$ws.Range('E24').Value = '  +0.84%  '
# Row 25
$ws.Range('D25').Value = '0.1552'
$ws.Range('E25').Value = '  -2.57%  '
# Row 26
$ws.Range('D26').Value = '8.953'
$ws.Range('E26').Value = '  -1.31%  '
# Row 27
$ws.Range('D27').Value = '159.75'
$ws.Range('E27').Value = '  -2.07%  '
# Row 28
$ws.Range('D28').Value = '18.57'
$ws.Range('E28').Value = '  +0.06%  '
# Row 29
$ws.Range('D29').Value = '1.503'
$ws.Range('E29').Value = '  -0.18%  '
# Row 30
$ws.Range('D30').Value = '4.264'
$ws.Range('E30').Value = '  -3.62%  '
# Row 31
$ws.Range('D31').Value = '4.249'
$ws.Range('E31').Value = '  -2.18%  '
# Row 32
$ws.Range('D32').Value = '1.227'
$ws.Range('E32').Value = '  +1.83%  '
# Row 33
$ws.Range('D33').Value = '0.05238'
$ws.Range('E33').Value = '  -2.53%  '
# Row 34
$ws.Range('D34').Value = '1.877'
$ws.Range('E34').Value = '  -3.85%  '
# Row 35
$ws.Range('D35').Value = '0.7422'
$ws.Range('E35').Value = '  -1.33%  '
# Row 36
$ws.Range('D36').Value = '1.156'
$ws.Range('E36').Value = '  -1.98%  '
# Row 37
$ws.Range('D37').Value = '2.721'
$ws.Range('E37').Value = '  +0.79%  '
# Row 38
$ws.Range('D38').Value = '0.01859'
$ws.Range('E38').Value = '  -1.13%  '
# Row 39
$ws.Range('D39').Value = '1.255.80'
$ws.Range('E39').Value = '  -1.97%  '
# Row 40
$ws.Range('D40').Value = '2.744'
$ws.Range('E40').Value = '  +0.12%  '
# Row 41
$ws.Range('D41').Value = '6.235'
$ws.Range('E41').Value = '  -5.28%  '
# Row 42
$ws.Range('D42').Value = '0.8977'
$ws.Range('E42').Value = '  +0.56%  '
# Row 43
$ws.Range('D43').Value = '110.62'
$ws.Range('E43').Value = '  +0.25%  '
# Row 44
$ws.Range('D44').Value = '70.96'
$ws.Range('E44').Value = '  -3.03%  '
# Row 45
$ws.Range('D45').Value = '1.006'
$ws.Range('E45').Value = '  +0.61%  '
# Row 46
$ws.Range('D46').Value = '2.033.76'
$ws.Range('E46').Value = '  +0.79%  '
# Row 47
$ws.Range('E47').Value = '  -4.57%  '
# Row 48
$ws.Range('D48').Value = '0.5215'
$ws.Range('E48').Value = '  +0.25%  '
# Row 49
$ws.Range('D49').Value = '1.773'
$ws.Range('E49').Value = '  -1.68%  '
# Row 50
$ws.Range('D50').Value = '9.271'
$ws.Range('E50').Value = '  -2.18%  '
# Row 51
$ws.Range('B51').Value = 'Frax'
$ws.Range('C51').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D51').Value = '1.007'
$ws.Range('E51').Value = '  +0.95%  '
